$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("F2").Value = 11.17
$ws.Range("E3").Value = 10.8
$ws.Range("E4").Value = 10.59
$ws.Range("F4").Value = 9.68
$ws.Range("G4").Value = 10.15
$ws.Range("C5").Value = 9.15
$ws.Range("D5").Value = 9.41
$ws.Range("F5").Value = 10.16
$ws.Range("B6").Value = 8.83
$ws.Range("D6").Value = 10.32
$ws.Range("E6").Value = 9.84
$ws.Range("H6").Value = 10.61
$ws.Range("D7").Value = 10
$ws.Range("H7").Value = 9.9
$ws.Range("F8").Value = 9.39
$ws.Range("G8").Value = 10.1
$ws.Range("J8").Value = 11
$ws.Range("H10").Value = 9
